# basic randomization of moves
# Randomize NumberOfBeats (column D) for a handful of moves in the
# dance-move spreadsheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D53").Value = 12   # cross break
$ws.Range("D55").Value = 16   # swivels
$ws.Range("D56").Value = 16   # McFreeBird

# Leave the sheet scrolled/selected where editing left off.
$ws.Range("D57").Select()
